$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 106; this shifts the existing rows
# 106-152 down to 108-154 (preserves all their data/styles as-is).
$ws.Rows("106:107").Insert()

# Populate the two newly inserted rows with the new weekly price entries
# (Vega Monumental Concepción, Frutilla, week of 2021-09-21 / serial 44460).

# Row 106: Calidad "Especial"
$ws.Cells.Item(106, 1).Value = 11
$ws.Cells.Item(106, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(106, 3).Value = "Bíobío"
$ws.Cells.Item(106, 4).Value = 44460
$ws.Cells.Item(106, 5).Value = 8
$ws.Cells.Item(106, 6).Value = "Fruta"
$ws.Cells.Item(106, 7).Value = 100101
$ws.Cells.Item(106, 8).Value = "Berries"
$ws.Cells.Item(106, 9).Value = 100112025
$ws.Cells.Item(106, 10).Value = "Frutilla"
$ws.Cells.Item(106, 11).Value = "Sin especificar"
$ws.Cells.Item(106, 12).Value = "Especial"
$ws.Cells.Item(106, 13).Value = 50
$ws.Cells.Item(106, 14).Value = 22000
$ws.Cells.Item(106, 15).Value = 22000
$ws.Cells.Item(106, 16).Value = 22000
$ws.Cells.Item(106, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(106, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(106, 19).Value = 3143
$ws.Cells.Item(106, 20).Value = 7

# Row 107: Calidad "Primera"
$ws.Cells.Item(107, 1).Value = 11
$ws.Cells.Item(107, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(107, 3).Value = "Bíobío"
$ws.Cells.Item(107, 4).Value = 44460
$ws.Cells.Item(107, 5).Value = 8
$ws.Cells.Item(107, 6).Value = "Fruta"
$ws.Cells.Item(107, 7).Value = 100101
$ws.Cells.Item(107, 8).Value = "Berries"
$ws.Cells.Item(107, 9).Value = 100112025
$ws.Cells.Item(107, 10).Value = "Frutilla"
$ws.Cells.Item(107, 11).Value = "Sin especificar"
$ws.Cells.Item(107, 12).Value = "Primera"
$ws.Cells.Item(107, 13).Value = 50
$ws.Cells.Item(107, 14).Value = 18000
$ws.Cells.Item(107, 15).Value = 18000
$ws.Cells.Item(107, 16).Value = 18000
$ws.Cells.Item(107, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(107, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(107, 19).Value = 2571
$ws.Cells.Item(107, 20).Value = 7
